# Apply small "want-to-go" count (column F) increments across sheets,
# matching the re-scraped data snapshot described by the commit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 373
$ws1.Range("F6").Value  = 844
$ws1.Range("F11").Value = 6237
$ws1.Range("F12").Value = 6237
$ws1.Range("F13").Value = 72
$ws1.Range("F14").Value = 471
$ws1.Range("F20").Value = 9377
$ws1.Range("F22").Value = 2531
$ws1.Range("F24").Value = 2341
$ws1.Range("F25").Value = 2505
$ws1.Range("F37").Value = 588
$ws1.Range("F42").Value = 1576

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 704

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 704
$ws4.Range("F7").Value  = 373
$ws4.Range("F10").Value = 844
$ws4.Range("F16").Value = 6237
$ws4.Range("F17").Value = 72
$ws4.Range("F21").Value = 9377
$ws4.Range("F24").Value = 2531
$ws4.Range("F26").Value = 2505
$ws4.Range("F37").Value = 588
